$d = $word.ActiveDocument

$replacements = @(
    @("636×4=", "237×3="),
    @("901×5=", "502×9="),
    @("245×7=", "147×6="),
    @("206×6=", "311×2="),
    @("771×2=", "646×7="),
    @("172×2=", "888×8="),
    @("144×8=", "177×7="),
    @("880×6=", "276×4="),
    @("527×3=", "126×3="),
    @("725×8=", "494×5="),
    @("707×9=", "365×7="),
    @("914×7=", "636×8="),
    @("729×4=", "138×8="),
    @("176×3=", "772×5="),
    @("220×5=", "308×8="),
    @("423×8=", "905×2="),
    @("151×6=", "367×6="),
    @("746×9=", "246×7="),
    @("632×5=", "413×3="),
    @("680×2=", "875×7="),
    @("666×3=", "596×9="),
    @("402×3=", "693×5="),
    @("392×6=", "311×9="),
    @("645×5=", "938×9="),
    @("197×8=", "618×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
